$d = $word.ActiveDocument
$d.Content.Find.Execute("Windows 8/10 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Windows 10 ", 2)
